$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking values so they are stored as text, not numbers
$textCells = @("D5","D6","D7","D8","D10","D11","D12","D16","D19","D20","D21","D22","D23","D24","D25","D30","D32","D33","D34","D35","D39","D40","D41","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '59.249.69'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '2.634.64'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '515.21'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").Value = '145.82'
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = '0.572'
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = '2.668.24'
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("D10").Value = '6.44'
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("D11").Value = '0.106'
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").Value = '0.339'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("D14").Value = '3.098.99'
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").Value = '59.232.59'
$ws.Range("E15").Value = '  -1.66%  '
$ws.Range("D16").Value = '21.21'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '2.657.23'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '4.60'
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").Value = '345.36'
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("D21").Value = '10.48'
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '61.43'
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("D26").Value = '2.760.11'
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").Value = '0.0₃0814'
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("D30").Value = '7.18'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").Value = '6.53'
$ws.Range("E32").Value = '  +9.10%  '
$ws.Range("D33").Value = '19.07'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").Value = '150.17'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("E36").Value = '  +12.93%  '
$ws.Range("E37").Value = '  +2.92%  '
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").Value = '0.864'
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").Value = '36.73'
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").Value = '3.72'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("D43").Value = '285.32'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").Value = '0.617'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.995'
$ws.Range("E45").Value = '  -0.19%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0986'
$ws.Range("E46").Value = '  -1.68%  '
$ws.Range("D47").Value = '19.65'
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").Value = '0.0542'
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("D49").Value = '0.0232'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = '4.74'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("E51").Value = '  -1.36%  '
